# daily auto push: 2025-10-05 13:28 UTC
# Append the new daily-ranking row (row 65) to the bottom of the data table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A holds a date-like string ("2025/10/05") that must stay plain text,
# matching the rest of the column (inline/shared string, not a date serial).
# Temporarily force a text format so Excel's auto-detection doesn't convert
# it to a date serial number, then clear the formatting again so the new
# cell doesn't end up with a style index unlike its siblings.
$ws.Range("A65").NumberFormat = "@"
$ws.Range("A65").Value = "2025/10/05"
$ws.Range("A65").ClearFormats()

$ws.Range("B65").Value = "日"
$ws.Range("C65").Value = 20
$ws.Range("D65").Value = 5
